$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-6 (false position / regula falsi method for f(x)=sin(x))
$data = @(
    @(1, [double]"-0.441888565985864", [double]"-2", [double]"1", [double]"-0.427647387672114", [double]"-0.909297426825682", [double]"0.8414709848078969", [double]"1.0005"),
    @(2, [double]"0.043976182274927", [double]"-0.441888565985864", [double]"1", [double]"0.0439620093551967", [double]"-0.427647387672114", [double]"0.8414709848078969", [double]"0.485864748260791"),
    @(3, [double]"-0.0013146680411116", [double]"-0.441888565985864", [double]"0.043976182274927", [double]"-0.0013146676624101", [double]"-0.427647387672114", [double]"0.0439620093551967", [double]"0.0452908503160387"),
    @(4, [double]"4.11161595739468e-07", [double]"-0.0013146680411116", [double]"0.043976182274927", [double]"4.11161595739457e-07", [double]"-0.0013146676624101", [double]"0.0439620093551967", [double]"0.0013150792027074"),
    @(5, [double]"-1.18401599288886e-13", [double]"-0.0013146680411116", [double]"4.11161595739468e-07", [double]"-1.18401599288886e-13", [double]"-0.0013146676624101", [double]"4.11161595739457e-07", [double]"4.11161714141068e-07")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}

# Remove old rows 7-11 (they no longer exist in the updated table)
$ws.Range("A7:H11").Delete() | Out-Null
